{"js": "// Update the reported scikit-learn / custom-OLS / gradient-descent / SGD\n// model coefficients that appear (twice, for the two duplicated sections,\n// plus once each for the other two) in the \"So the model parameters: ...\"\n// paragraphs.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"1.0075 and coefficients\", \"1.0132 and coefficients\"],\n  [\"2.981, and\", \"2.9918, and\"],\n  [\"1.9891\", \"1.9769\"],\n  [\"1.0079 and coefficients\", \"1.0125 and coefficients\"],\n  [\"2.9795, and\", \"2.9906, and\"],\n  [\"1.9885\", \"1.9762\"],\n  [\"array([1.05822114]) and coefficients\", \"array([1.01726883]) and coefficients\"],\n  [\"array([2.98804002]), and\", \"array([2.99623763]), and\"],\n  [\"array([2.02031678])\", \"array([1.94734374])\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the reported scikit-learn / custom-OLS / gradient-descent / SGD\n# model coefficients that appear (twice, for the two duplicated sections,\n# plus once each for the other two) in the \"So the model parameters: ...\"\n# paragraphs.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"1.0075 and coefficients\"; New = \"1.0132 and coefficients\" },\n    @{ Old = \"2.981, and\"; New = \"2.9918, and\" },\n    @{ Old = \"1.9891\"; New = \"1.9769\" },\n    @{ Old = \"1.0079 and coefficients\"; New = \"1.0125 and coefficients\" },\n    @{ Old = \"2.9795, and\"; New = \"2.9906, and\" },\n    @{ Old = \"1.9885\"; New = \"1.9762\" },\n    @{ Old = \"array([1.05822114]) and coefficients\"; New = \"array([1.01726883]) and coefficients\" },\n    @{ Old = \"array([2.98804002]), and\"; New = \"array([2.99623763]), and\" },\n    @{ Old = \"array([2.02031678])\"; New = \"array([1.94734374])\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute(\n        $r.Old,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $r.New,\n        2\n    )\n}\n"}
